$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Force the cell to store $val as literal text, matching the
    # workbook's existing convention (numeric-looking values kept as
    # strings), while preserving the cell's original (General) style.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value2 = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "277.01"
Set-TextValue "D3" "21.18"
Set-TextValue "D4" "6.272"
Set-TextValue "D5" "0.06206"
Set-TextValue "D6" "3.555"
Set-TextValue "D7" "1.537"
Set-TextValue "D8" "6.579"
Set-TextValue "D9" "0.8279"
Set-TextValue "D10" "0.1667"
Set-TextValue "D11" "0.08294"
Set-TextValue "D12" "0.03507"
Set-TextValue "D13" "0.03197"
Set-TextValue "D14" "0.09159"
Set-TextValue "D15" "3.765"
Set-TextValue "D16" "0.001636"
Set-TextValue "D17" "0.04696"
Set-TextValue "D18" "0.006286"
Set-TextValue "D19" "0.006220"
Set-TextValue "D22" "3.721"
Set-TextValue "D23" "2.313"
Set-TextValue "D24" "0.01397"
Set-TextValue "D25" "0.3290"
Set-TextValue "D28" "0.0002735"
Set-TextValue "D40" "0.04746"
$ws.Range("B41").Value = "CEJI"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D41" "0.005197"
$ws.Range("E41").Value = "40CEJICEJI"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D42" "0.007085"
$ws.Range("E42").Value = "41KickTokenKICK"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1120"
$ws.Range("E43").Value = "42BKEXTokenBKK"
Set-TextValue "D44" "0.01140"
Set-TextValue "D45" "0.00006379"
Set-TextValue "D47" "0.7226"
Set-TextValue "D48" "0.001400"
Set-TextValue "D49" "0.00001899"
Set-TextValue "D50" "0.01239"
